# Applies FlashScore odds update for 2024-10-30 workbook
# - Row 2: updates Odd_Over25_FT (Q2) and Odd_Under25_FT (R2)
# - Rows 5-7: match list shifted (one fixture removed, one added) with refreshed odds
# - Row 8: refreshed odds for existing fixture

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: Odd_Over25_FT / Odd_Under25_FT refresh ---
$ws.Range("Q2").Value = 1.84
$ws.Range("R2").Value = 2.06

# --- Row 5: full fixture refresh ---
$row5 = New-Object 'object[,]' 1,56
$row5[0,0] = "MyCLs2Y8"
$row5[0,1] = "30/10/2024"
$row5[0,2] = "14:00"
$row5[0,3] = "EGYPT - PREMIER LEAGUE"
$row5[0,4] = "El Gaish"
$row5[0,5] = "Al Masry"
$row5[0,6] = 4.35
$row5[0,7] = 2.72
$row5[0,8] = 2.07
$row5[0,9] = 5.1
$row5[0,10] = 1.83
$row5[0,11] = 2.72
$row5[0,12] = 1.15
$row5[0,13] = 4.7
$row5[0,14] = 1.65
$row5[0,15] = 2.15
$row5[0,16] = 2.82
$row5[0,17] = 1.38
$row5[0,18] = 1.62
$row5[0,19] = 2.22
$row5[0,20] = 2.3
$row5[0,21] = 1.55
$row5[0,22] = 7.9
$row5[0,23] = 22
$row5[0,24] = 16
$row5[0,25] = 90
$row5[0,26] = 65
$row5[0,27] = 90
$row5[0,28] = 4.7
$row5[0,29] = 5.7
$row5[0,30] = 21
$row5[0,31] = 150
$row5[0,32] = 201
$row5[0,33] = 5
$row5[0,34] = 8.25
$row5[0,35] = 9.25
$row5[0,36] = 19.5
$row5[0,37] = 22
$row5[0,38] = 45
$row5[0,39] = 5.9
$row5[0,40] = 30
$row5[0,41] = 40
$row5[0,42] = 200
$row5[0,43] = 300
$row5[0,44] = 500
$row5[0,45] = 2.18
$row5[0,46] = 8
$row5[0,47] = 100
$row5[0,48] = 3.65
$row5[0,49] = 11
$row5[0,50] = 24
$row5[0,51] = 45
$row5[0,52] = 100
$row5[0,53] = 400
$row5[0,54] = 51
$row5[0,55] = 51
$ws.Range("A5:BD5").Value = $row5

# --- Row 6: full fixture refresh ---
$row6 = New-Object 'object[,]' 1,56
$row6[0,0] = "lrf80kOO"
$row6[0,1] = "30/10/2024"
$row6[0,2] = "13:00"
$row6[0,3] = "FINLAND - VEIKKAUSLIIGA"
$row6[0,4] = "Haka"
$row6[0,5] = "SJK"
$row6[0,6] = 3.6
$row6[0,7] = 3.5
$row6[0,8] = 1.95
$row6[0,9] = 4
$row6[0,10] = 2.3
$row6[0,11] = 2.6
$row6[0,12] = 1.04
$row6[0,13] = 13
$row6[0,14] = 1.2
$row6[0,15] = 4.33
$row6[0,16] = 1.7
$row6[0,17] = 2.1
$row6[0,18] = 1.33
$row6[0,19] = 3.25
$row6[0,20] = 1.62
$row6[0,21] = 2.2
$row6[0,22] = 13
$row6[0,23] = 21
$row6[0,24] = 13
$row6[0,25] = 41
$row6[0,26] = 26
$row6[0,27] = 29
$row6[0,28] = 13
$row6[0,29] = 7
$row6[0,30] = 12
$row6[0,31] = 41
$row6[0,32] = 151
$row6[0,33] = 9
$row6[0,34] = 11
$row6[0,35] = 8.5
$row6[0,36] = 17
$row6[0,37] = 15
$row6[0,38] = 21
$row6[0,39] = 5.5
$row6[0,40] = 19
$row6[0,41] = 23
$row6[0,42] = 51
$row6[0,43] = 67
$row6[0,44] = 151
$row6[0,45] = 3.25
$row6[0,46] = 7.5
$row6[0,47] = 41
$row6[0,48] = 4.33
$row6[0,49] = 10
$row6[0,50] = 19
$row6[0,51] = 34
$row6[0,52] = 51
$row6[0,53] = 101
$row6[0,54] = 451
$row6[0,55] = 81
$ws.Range("A6:BD6").Value = $row6

# --- Row 7: full fixture refresh ---
$row7 = New-Object 'object[,]' 1,56
$row7[0,0] = "CWHk9jo1"
$row7[0,1] = "30/10/2024"
$row7[0,2] = "15:30"
$row7[0,3] = "ITALY - SERIE B"
$row7[0,4] = "Sudtirol"
$row7[0,5] = "Frosinone"
$row7[0,6] = 2.55
$row7[0,7] = 2.75
$row7[0,8] = 3.1
$row7[0,9] = 3.5
$row7[0,10] = 1.83
$row7[0,11] = 4
$row7[0,12] = 1.14
$row7[0,13] = 5.5
$row7[0,14] = 1.57
$row7[0,15] = 2.25
$row7[0,16] = 2.88
$row7[0,17] = 1.4
$row7[0,18] = 1.62
$row7[0,19] = 2.2
$row7[0,20] = 2.2
$row7[0,21] = 1.62
$row7[0,22] = 6
$row7[0,23] = 11
$row7[0,24] = 11
$row7[0,25] = 26
$row7[0,26] = 26
$row7[0,27] = 41
$row7[0,28] = 5.5
$row7[0,29] = 5.5
$row7[0,30] = 19
$row7[0,31] = 81
$row7[0,32] = 501
$row7[0,33] = 7
$row7[0,34] = 13
$row7[0,35] = 13
$row7[0,36] = 34
$row7[0,37] = 34
$row7[0,38] = 51
$row7[0,39] = 4.33
$row7[0,40] = 17
$row7[0,41] = 34
$row7[0,42] = 51
$row7[0,43] = 101
$row7[0,44] = 351
$row7[0,45] = 2.2
$row7[0,46] = 9.5
$row7[0,47] = 81
$row7[0,48] = 4.75
$row7[0,49] = 19
$row7[0,50] = 34
$row7[0,51] = 67
$row7[0,52] = 126
$row7[0,53] = 401
$row7[0,54] = 81
$row7[0,55] = 81
$ws.Range("A7:BD7").Value = $row7

# --- Row 8: odds refresh ---
$ws.Range("J8").Value = 5.8
$ws.Range("K8").Value = 2.55
$ws.Range("N8").Value = 10.5
$ws.Range("O8").Value = 1.15
$ws.Range("P8").Value = 5
$ws.Range("Q8").Value = 1.47
$ws.Range("R8").Value = 2.55
$ws.Range("S8").Value = 1.27
$ws.Range("T8").Value = 3.55
$ws.Range("U8").Value = 1.65
$ws.Range("V8").Value = 2.1
$ws.Range("W8").Value = 21
$ws.Range("Y8").Value = 22
$ws.Range("AB8").Value = 55
$ws.Range("AC8").Value = 10.5
$ws.Range("AD8").Value = 10.5
$ws.Range("AE8").Value = 17.5
$ws.Range("AF8").Value = 65
$ws.Range("AG8").Value = 400
$ws.Range("AH8").Value = 8.5
$ws.Range("AK8").Value = 10.75
$ws.Range("AM8").Value = 22
$ws.Range("AR8").Value = 175
$ws.Range("AS8").Value = 300
$ws.Range("AT8").Value = 3.55
$ws.Range("AU8").Value = 7.5
$ws.Range("AW8").Value = 3.5
$ws.Range("AY8").Value = 13.5
$ws.Range("BA8").Value = 35

